$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by
# the type inference on Range.Value are forced to Text first, written,
# then ClearFormats() restores the default (General) cell style so no
# stray number format is left behind.

$ws.Range("D2").Value = "41.143.95"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "2.429.74"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.46"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "89.48"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.68%  "
$ws.Range("E7").Value = "  -1.97%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.498"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0837"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.07"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.26%  "
$ws.Range("E12").Value = "  -1.44%  "
$ws.Range("D13").Value = "2.805.58"
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.76"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").Value = "2.418.87"
$ws.Range("E16").Value = "  -3.04%  "
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("D18").Value = "41.093.60"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").Value = "0.0₃0927"
$ws.Range("E19").Value = "  -1.86%  "
$ws.Range("E20").Value = "  -2.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.08"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.09"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.12"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.54%  "
$ws.Range("E24").Value = "  -1.34%  "
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("E26").Value = "  -2.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.19"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.40%  "
$ws.Range("E28").Value = "  -1.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.62"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.63"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.81"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.27"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.99%  "
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0746"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.50"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.96"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.25%  "
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("E39").Value = "  -2.08%  "
$ws.Range("E40").Value = "  -2.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.90"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.88%  "
$ws.Range("D42").Value = "1.995.83"
$ws.Range("E42").Value = "  +0.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.27"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -7.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.59"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.49%  "
$ws.Range("E45").Value = "  -2.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.90"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.56"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.40%  "
$ws.Range("D48").Value = "2.666.99"
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "94.82"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.56"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.82"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.52%  "
